$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.876.87"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "3.390.17"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.28%  "
$ws.Range("E7").Value = "  +3.94%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.386.20"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "3.976.91"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("D16").Value = "65.938.40"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "3.399.21"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.526"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.857"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.82%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").Value = "2.668.07"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0677"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "331.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.18%  "
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.43%  "
